$wb = $excel.ActiveWorkbook

# --- Add Denmark sheet (copied from UK) ---
$uk = $wb.Worksheets.Item("UK")
$uk.Copy($null, $wb.Worksheets.Item($wb.Worksheets.Count)) | Out-Null
$denmark = $wb.Worksheets.Item($wb.Worksheets.Count)
$denmark.Name = "Denmark"
$denmark.Range("B2").Value = "Denmark market"
$denmark.Range("B4").Value = "NGC-3446/T2009"
$denmark.Range("A1:XFD1048576").Select() | Out-Null

# --- Add Sweden sheet (copied from UK) ---
$uk.Copy($null, $wb.Worksheets.Item($wb.Worksheets.Count)) | Out-Null
$sweden = $wb.Worksheets.Item($wb.Worksheets.Count)
$sweden.Name = "Sweden"
$sweden.Range("B2").Value = "Sweden Market"
$sweden.Range("B4").ClearFormats() | Out-Null
$sweden.Range("B4").Value = "NGC-3465/T2021"
$sweden.Range("A1:XFD1048576").Select() | Out-Null

# --- Add Norway sheet (copied from UK) ---
$uk.Copy($null, $wb.Worksheets.Item($wb.Worksheets.Count)) | Out-Null
$norway = $wb.Worksheets.Item($wb.Worksheets.Count)
$norway.Name = "Norway"
$norway.Range("B4").Value = "NGC-3464/T1924"
$norway.Range("B2").Value = "Norway Market"
$norway.Range("B7").Select() | Out-Null

Write-Host "done"
